# Zombono 0.0.11 version-scheduling sheet update.
#
# Adds two new task rows to the schedule table:
#   - "Make Master servers work" (Feature, Netservuces)
#   - "cl_console_line_length"   (Feature)
# Increased console line length from 38 to 128 characters
# (this will be a CVar in the future).
#
# Rows below the table already reserved the "D" column formatting on
# specific row numbers (10/12/15/16/17), so rather than inserting whole
# rows (which would drag that per-row formatting along with it), the
# table is rewritten in place directly to its final A:B contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value2  = "Make Master servers work"
$ws.Range("B9").Value2  = "Feature, Netservuces"

$ws.Range("A10").Value2 = "BrowseServersUI"
$ws.Range("B10").Value2 = "Feature, Netservices"

$ws.Range("A11").Value2 = "Text Engine - scaled coordinates"
$ws.Range("B11").Value2 = "Feature"

$ws.Range("A12").Value2 = "Start Waves mode programming"
$ws.Range("B12").Value2 = "Feature"

$ws.Range("A13").Value2 = "Properly split out client.h, server.h"
$ws.Range("B13").Value2 = "Refactoring"

$ws.Range("A14").Value2 = "Allow people to see what team a player is"
$ws.Range("B14").Value2 = "Feature"

$ws.Range("A15").Value2 = "cl_console_line_length"
$ws.Range("B15").Value2 = "Feature"

$ws.Range("A16").Value2 = "z_waves_port working (THEY COME FROM THE SEA/!?!?!?!?!) - z_tdm_spire finished"
$ws.Range("B16").Value2 = "Content"

$ws.Range("A17").Value2 = "Finish z_warehouse easter egg"
$ws.Range("B17").Value2 = "Content"

# Match the recorded selection left behind by the edit.
$ws.Range("B15").Select()
